$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 440.66666
$ws.Range("I5").Value = 292.55554
$ws.Range("K5").Value = 292.55554
$ws.Range("M5").Value = -177.55554

$ws.Range("H17").Value = 1053592.2
$ws.Range("J17").Value = 1053592.2
$ws.Range("L17").Value = 3160776.6
$ws.Range("N17").Value = -3161112.6

$ws.Range("H38").Value = 415.8
$ws.Range("I38").Value = 270
$ws.Range("J38").Value = 999
$ws.Range("K38").Value = 810
$ws.Range("L38").Value = 2997
$ws.Range("M38").Value = -438
$ws.Range("N38").Value = -3741

$ws.Range("H40").Value = 4074.875
$ws.Range("J40").Value = 4433.1665
$ws.Range("L40").Value = 4433.1665
$ws.Range("N40").Value = -4783.1665

$ws.Range("H42").Value = 927.1429000000001
$ws.Range("J42").Value = 1199
$ws.Range("L42").Value = 3597
$ws.Range("N42").Value = -4057

$ws.Range("H96").Value = 166668000
$ws.Range("I96").Value = 1949.75
$ws.Range("K96").Value = 5849.25
$ws.Range("M96").Value = -4476.25

$ws.Range("H98").Value = 4562.375
$ws.Range("I98").Value = 4333.1665
$ws.Range("J98").Value = 5250
$ws.Range("K98").Value = 4333.1665
$ws.Range("L98").Value = 5250
$ws.Range("M98").Value = -2835.1665
$ws.Range("N98").Value = -8246

$ws.Range("H99").Value = 111114616
$ws.Range("J99").Value = 333343140
$ws.Range("L99").Value = 1000029420
$ws.Range("N99").Value = -1000032416

$ws.Range("H103").Value = 384.69232
$ws.Range("J103").Value = 579.3333
$ws.Range("L103").Value = 1737.9999
$ws.Range("N103").Value = -2909.9999

$ws.Range("H107").Value = 2770.2173
$ws.Range("I107").Value = 1242.3684
$ws.Range("K107").Value = 1242.3684
$ws.Range("M107").Value = 677.6315999999999

$ws.Range("H115").Value = 1957.3
$ws.Range("I115").Value = 1063.8889
$ws.Range("K115").Value = 3191.6667
$ws.Range("M115").Value = -1624.6667

$ws.Range("H122").Value = 4562.375
$ws.Range("I122").Value = 4333.1665
$ws.Range("J122").Value = 5250
$ws.Range("K122").Value = 12999.4995
$ws.Range("L122").Value = 15750
$ws.Range("M122").Value = -10549.4995
$ws.Range("N122").Value = -20650

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3095.8918
$ws.Range("I61").Value = 2804.9375
$ws.Range("J61").Value = 4958
$ws.Range("K61").Value = 2804.9375
$ws.Range("L61").Value = 4958
$ws.Range("M61").Value = -2592.9375
$ws.Range("N61").Value = -5382

$ws.Range("H63").Value = 2106.5
$ws.Range("I63").Value = 2033.25
$ws.Range("J63").Value = 2399.5
$ws.Range("K63").Value = 2033.25
$ws.Range("L63").Value = 2399.5
$ws.Range("M63").Value = -1347.25
$ws.Range("N63").Value = -3771.5

$ws.Range("H66").Value = 2106.5
$ws.Range("I66").Value = 2033.25
$ws.Range("J66").Value = 2399.5
$ws.Range("K66").Value = 10166.25
$ws.Range("L66").Value = 11997.5
$ws.Range("M66").Value = -6734.25
$ws.Range("N66").Value = -18861.5

$ws.Range("H136").Value = 3095.8918
$ws.Range("I136").Value = 2804.9375
$ws.Range("J136").Value = 4958
$ws.Range("K136").Value = 8414.8125
$ws.Range("L136").Value = 14874
$ws.Range("M136").Value = -5864.8125
$ws.Range("N136").Value = -19974

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 12361.546
$ws.Range("I94").Value = 5997.5
$ws.Range("J94").Value = 19998.4
$ws.Range("K94").Value = 5997.5
$ws.Range("L94").Value = 19998.4
$ws.Range("M94").Value = -5546.5
$ws.Range("N94").Value = -20900.4

$ws.Range("H141").Value = 207200
$ws.Range("J141").Value = 207200
$ws.Range("L141").Value = 207200
$ws.Range("N141").Value = -217560

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4469
$ws.Range("I31").Value = 1709.5
$ws.Range("J31").Value = 7016.231
$ws.Range("K31").Value = 1709.5
$ws.Range("L31").Value = 7016.231
$ws.Range("M31").Value = -1414.5
$ws.Range("N31").Value = -7606.231

$ws.Range("H34").Value = 4469
$ws.Range("I34").Value = 1709.5
$ws.Range("J34").Value = 7016.231
$ws.Range("K34").Value = 1709.5
$ws.Range("L34").Value = 7016.231
$ws.Range("M34").Value = -1507.5
$ws.Range("N34").Value = -7420.231

$ws.Range("H134").Value = 56057420
$ws.Range("I134").Value = 63061696
$ws.Range("J134").Value = 23208
$ws.Range("K134").Value = 189185088
$ws.Range("L134").Value = 69624
$ws.Range("M134").Value = -189182553
$ws.Range("N134").Value = -74694

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1899.5
$ws.Range("I17").Value = 1899.5
$ws.Range("K17").Value = 5698.5
$ws.Range("M17").Value = -5529.5

$ws.Range("H25").Value = 1833
$ws.Range("I25").Value = 1749.5
$ws.Range("K25").Value = 5248.5
$ws.Range("M25").Value = -5079.5

$ws.Range("H30").Value = 1833
$ws.Range("I30").Value = 1749.5
$ws.Range("K30").Value = 5248.5
$ws.Range("M30").Value = -5146.5

$ws.Range("H37").Value = 135996
$ws.Range("J37").Value = 135996
$ws.Range("L37").Value = 407988
$ws.Range("N37").Value = -408212

$ws.Range("H39").Value = 3269
$ws.Range("J39").Value = 3269
$ws.Range("L39").Value = 9807
$ws.Range("N39").Value = -10395

$ws.Range("H131").Value = 27958752
$ws.Range("J131").Value = 25644036
$ws.Range("L131").Value = 76932108
$ws.Range("N131").Value = -76942188

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4270.5713
$ws.Range("I93").Value = 4542.2856
$ws.Range("K93").Value = 4542.2856
$ws.Range("M93").Value = -3294.2856

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 26496.666
$ws.Range("I62").Value = 32750
$ws.Range("J62").Value = 23370
$ws.Range("K62").Value = 32750
$ws.Range("L62").Value = 23370
$ws.Range("M62").Value = -32126
$ws.Range("N62").Value = -24618

$ws.Range("H65").Value = 26496.666
$ws.Range("I65").Value = 32750
$ws.Range("J65").Value = 23370
$ws.Range("K65").Value = 163750
$ws.Range("L65").Value = 116850
$ws.Range("M65").Value = -160630
$ws.Range("N65").Value = -123090

$ws.Range("H132").Value = 6581
$ws.Range("I132").Value = 3137.394
$ws.Range("J132").Value = 34990.75
$ws.Range("K132").Value = 9412.181999999999
$ws.Range("L132").Value = 104972.25
$ws.Range("M132").Value = -6882.181999999999
$ws.Range("N132").Value = -110032.25

